$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly scoreboard rows (week 8), appended after the existing last row (323).
$rows = @(
    @{ Row=324; A="Matt";   B=45506; C="Ride"; D=78;  E=21.26; F=1207; G=3;   H=29; I=35; J=8; K=1; L="Agile Antelope";  M=8 },
    @{ Row=325; A="Steven"; B=45506; C="Walk"; D=28;  E=1.27;  F=56;   G=28;  H=0;  I=0;  J=0; K=0; L="Brave Leopard";   M=8 },
    @{ Row=326; A="Steven"; B=45506; C="Walk"; D=82;  E=3.63;  F=82;   G=82;  H=0;  I=0;  J=0; K=0; L="Brave Leopard";   M=8 },
    @{ Row=327; A="Matt";   B=45507; C="Run";  D=45;  E=4.58;  F=207;  G=1;   H=23; I=10; J=8; K=0; L="Agile Antelope";  M=8 },
    @{ Row=328; A="Matt";   B=45507; C="Walk"; D=3;   E=0.12;  F=13;   G=3;   H=0;  I=0;  J=0; K=0; L="Agile Antelope";  M=8 },
    @{ Row=329; A="Steven"; B=45507; C="Walk"; D=142; E=6;     F=472;  G=142; H=0;  I=0;  J=0; K=0; L="Brave Leopard";   M=8 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column B carries the existing date style (numFmtId 14 / m/d/yyyy); copy the
    # format from the last pre-existing date cell instead of letting Excel mint a
    # brand-new style entry for it.
    $ws.Range("B323").Copy()
    $ws.Cells.Item($rowNum, 2).PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value2 = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
}

# Move the active selection to the new first blank row below the appended
# data and nudge the frozen pane's scroll position down to match, matching
# the saved view state (header row stays frozen via ySplit=1).
$win = $excel.ActiveWindow
$ws.Range("A330").Select() | Out-Null
$win.ScrollRow = 301
$null = $win.ScrollRow
